$wb = $excel.ActiveWorkbook

# --- Costs and Revenues ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76023.24100800102
$ws.Range("C2").Value = 32200
$ws.Range("D2").Value = 9272.289645445984
$ws.Range("E2").Value = 2370
$ws.Range("F2").Value = 9198.629642145857

# --- Capacities ---
$ws = $wb.Worksheets.Item("Capacities")
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 103
$ws.Range("B4").Value = 161
$ws.Range("C4").Value = 161
$ws.Range("D4").Value = 161

# --- PV Dispatch ---
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("G2").Value = 20.6
$ws.Range("H2").Value = 41.2
$ws.Range("I2").Value = 51.5
$ws.Range("J2").Value = 61.8
$ws.Range("K2").Value = 72.09999999999999
$ws.Range("L2").Value = 82.40000000000001
$ws.Range("M2").Value = 92.7
$ws.Range("N2").Value = 103
$ws.Range("O2").Value = 92.7
$ws.Range("P2").Value = 82.40000000000001
$ws.Range("Q2").Value = 72.09999999999999
$ws.Range("R2").Value = 51.5
$ws.Range("S2").Value = 30.9
$ws.Range("T2").Value = 20.6
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 82.40000000000001
$ws.Range("L3").Value = 92.7
$ws.Range("M3").Value = 103
$ws.Range("N3").Value = 26
$ws.Range("O3").Value = 25.657177838996
$ws.Range("P3").Value = 51.5
$ws.Range("Q3").Value = 51.5
$ws.Range("R3").Value = 30.9
$ws.Range("S3").Value = 20.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 41.2
$ws.Range("M4").Value = 31.88312417100287
$ws.Range("N4").Value = 82.40000000000001
$ws.Range("O4").Value = 72.09999999999999
$ws.Range("P4").Value = 41.2

# --- Battery Input ---
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("G2").Value = 12.8
$ws.Range("H2").Value = 28.2
$ws.Range("I2").Value = 20.3
$ws.Range("J2").Value = 22.8
$ws.Range("K2").Value = 46.1
$ws.Range("L2").Value = 61.6
$ws.Range("M2").Value = 69.3
$ws.Range("N2").Value = 77
$ws.Range("O2").Value = 61.5
$ws.Range("P2").Value = 53.8
$ws.Range("Q2").Value = 46.1
$ws.Range("R2").Value = 17.7
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 82.40000000000001
$ws.Range("L3").Value = 92.7
$ws.Range("M3").Value = 79.59999999999999
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 25.657177838996
$ws.Range("P3").Value = 22.9
$ws.Range("Q3").Value = 25.5
$ws.Range("R3").Value = 30.9
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 41.2
$ws.Range("M4").Value = 8.483124171002865
$ws.Range("N4").Value = 82.40000000000001
$ws.Range("O4").Value = 72.09999999999999
$ws.Range("P4").Value = 41.2

# --- Battery Output ---
$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 10.7
$ws.Range("T2").Value = 21.70772000002855
$ws.Range("S3").Value = 21

# --- State of Charge ---
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("B2").Value = 187.8909090909091
$ws.Range("C2").Value = 168.1939393939394
$ws.Range("D2").Value = 155.0626262626263
$ws.Range("E2").Value = 141.9313131313131
$ws.Range("F2").Value = 128.8
$ws.Range("G2").Value = 141.472
$ws.Range("H2").Value = 169.39
$ws.Range("I2").Value = 189.487
$ws.Range("J2").Value = 212.059
$ws.Range("K2").Value = 257.698
$ws.Range("L2").Value = 318.682
$ws.Range("M2").Value = 387.289
$ws.Range("N2").Value = 463.519
$ws.Range("O2").Value = 524.404
$ws.Range("P2").Value = 577.6660000000288
$ws.Range("Q2").Value = 623.3050000000288
$ws.Range("R2").Value = 640.8280000000289
$ws.Range("S2").Value = 630.019919191948
$ws.Range("T2").Value = 608.0929292929293
$ws.Range("U2").Value = 489.9111111111111
$ws.Range("V2").Value = 391.4262626262627
$ws.Range("W2").Value = 312.6383838383838
$ws.Range("X2").Value = 260.1131313131313
$ws.Range("Y2").Value = 220.7191919191919
$ws.Range("B3").Value = 181.3252525252525
$ws.Range("C3").Value = 161.6282828282828
$ws.Range("D3").Value = 148.4969696969697
$ws.Range("E3").Value = 148.4969696969697
$ws.Range("F3").Value = 148.4969696969697
$ws.Range("G3").Value = 128.8
$ws.Range("H3").Value = 128.8
$ws.Range("I3").Value = 128.8
$ws.Range("J3").Value = 128.8
$ws.Range("K3").Value = 210.376
$ws.Range("L3").Value = 302.149
$ws.Range("M3").Value = 380.953
$ws.Range("N3").Value = 380.953
$ws.Range("O3").Value = 406.353606060606
$ws.Range("P3").Value = 429.024606060606
$ws.Range("Q3").Value = 454.269606060606
$ws.Range("R3").Value = 484.860606060606
$ws.Range("S3").Value = 463.6484848484848
$ws.Range("T3").Value = 332.3353535353535
$ws.Range("U3").Value = 332.3353535353535
$ws.Range("V3").Value = 332.3353535353535
$ws.Range("W3").Value = 253.5474747474748
$ws.Range("X3").Value = 253.5474747474748
$ws.Range("Y3").Value = 214.1535353535353
$ws.Range("B4").Value = 168.1939393939394
$ws.Range("C4").Value = 148.4969696969697
$ws.Range("D4").Value = 148.4969696969697
$ws.Range("E4").Value = 148.4969696969697
$ws.Range("F4").Value = 148.4969696969697
$ws.Range("G4").Value = 128.8
$ws.Range("H4").Value = 128.8
$ws.Range("I4").Value = 128.8
$ws.Range("J4").Value = 128.8
$ws.Range("K4").Value = 169.588
$ws.Range("L4").Value = 169.588
$ws.Range("M4").Value = 177.9862929292929
$ws.Range("N4").Value = 259.5622929292929
$ws.Range("O4").Value = 330.9412929292929
$ws.Range("P4").Value = 371.7292929292929
$ws.Range("Q4").Value = 371.7292929292929
$ws.Range("R4").Value = 371.7292929292929
$ws.Range("S4").Value = 371.7292929292929
$ws.Range("T4").Value = 240.4161616161616
$ws.Range("U4").Value = 240.4161616161616
$ws.Range("V4").Value = 240.4161616161616
$ws.Range("W4").Value = 240.4161616161616
$ws.Range("X4").Value = 240.4161616161616
$ws.Range("Y4").Value = 201.0222222222222

# --- Feed in from Type 2 ---
$ws = $wb.Worksheets.Item("Feed in from Type 2")
$ws.Range("T2").Value = 0

# --- Feed in from Type 3 ---
$ws = $wb.Worksheets.Item("Feed in from Type 3")
$ws.Range("T2").Value = 9.692279999971451

Write-Host "Edit complete"